$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.672.16"
$ws.Range("E2").Value = "  +3.10%  "

# Row 3
$ws.Range("D3").Value = "2.204.34"
$ws.Range("E3").Value = "  +2.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'253.04"
$ws.Range("E5").Value = "  +6.67%  "

# Row 6
$ws.Range("D6").Value = "'0.614"
$ws.Range("E6").Value = "  +1.48%  "

# Row 7
$ws.Range("D7").Value = "'74.48"
$ws.Range("E7").Value = "  +4.31%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.589"
$ws.Range("E9").Value = "  +2.90%  "

# Row 10
$ws.Range("D10").Value = "'40.36"
$ws.Range("E10").Value = "  +2.24%  "

# Row 11
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = "  +2.49%  "

# Row 12
$ws.Range("D12").Value = "'6.86"
$ws.Range("E12").Value = "  +3.11%  "

# Row 13
$ws.Range("E13").Value = "  +1.79%  "

# Row 14
$ws.Range("D14").Value = "2.529.73"
$ws.Range("E14").Value = "  +1.99%  "

# Row 15
$ws.Range("D15").Value = "'14.38"
$ws.Range("E15").Value = "  +2.45%  "

# Row 16
$ws.Range("D16").Value = "2.216.31"
$ws.Range("E16").Value = "  +3.49%  "

# Row 17
$ws.Range("D17").Value = "'0.777"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "42.557.62"
$ws.Range("E18").Value = "  +3.24%  "

# Row 19
$ws.Range("D19").Value = "'0.0000103"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20
$ws.Range("D20").Value = "'71.14"
$ws.Range("E20").Value = "  +2.49%  "

# Row 21
$ws.Range("D21").Value = "'5.94"
$ws.Range("E21").Value = "  +3.35%  "

# Row 22
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'2.19"
$ws.Range("E22").Value = "  +10.48%  "

# Row 23
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'9.62"
$ws.Range("E23").Value = "  -1.59%  "

# Row 24
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'227.15"
$ws.Range("E24").Value = "  +0.12%  "

# Row 25
$ws.Range("E25").Value = "  -0.14%  "

# Row 26
$ws.Range("D26").Value = "'10.71"
$ws.Range("E26").Value = "  +0.96%  "

# Row 27
$ws.Range("E27").Value = "  +2.65%  "

# Row 28
$ws.Range("E28").Value = "  +2.82%  "

# Row 29
$ws.Range("E29").Value = "  +1.41%  "

# Row 30
$ws.Range("D30").Value = "'37.73"
$ws.Range("E30").Value = "  +16.34%  "

# Row 31
$ws.Range("E31").Value = "  -0.93%  "

# Row 32
$ws.Range("D32").Value = "'20.12"
$ws.Range("E32").Value = "  +4.18%  "

# Row 33
$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  +4.78%  "

# Row 34
$ws.Range("D34").Value = "'5.20"
$ws.Range("E34").Value = "  +2.12%  "

# Row 35
$ws.Range("D35").Value = "'0.121"
$ws.Range("E35").Value = "  +0.38%  "

# Row 36
$ws.Range("D36").Value = "'0.108"
$ws.Range("E36").Value = "  +4.61%  "

# Row 37
$ws.Range("D37").Value = "'4.38"
$ws.Range("E37").Value = "  +3.95%  "

# Row 38
$ws.Range("D38").Value = "'0.0332"
$ws.Range("E38").Value = "  +12.26%  "

# Row 39
$ws.Range("D39").Value = "'12.17"
$ws.Range("E39").Value = "  +2.43%  "

# Row 40
$ws.Range("D40").Value = "'2.08"
$ws.Range("E40").Value = "  +1.14%  "

# Row 41
$ws.Range("D41").Value = "'5.27"
$ws.Range("E41").Value = "  -0.56%  "

# Row 42
$ws.Range("D42").Value = "'0.198"
$ws.Range("E42").Value = "  +5.76%  "

# Row 43
$ws.Range("B43").Value = "WOONetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D43").Value = "'0.490"
$ws.Range("E43").Value = "  +27.90%  "

# Row 44
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'59.25"
$ws.Range("E44").Value = "  +1.05%  "

# Row 45
$ws.Range("D45").Value = "'103.28"
$ws.Range("E45").Value = "  +8.11%  "

# Row 46
$ws.Range("D46").Value = "'8.37"
$ws.Range("E46").Value = "  +0.48%  "

# Row 47
$ws.Range("D47").Value = "'0.0981"
$ws.Range("E47").Value = "  +3.15%  "

# Row 48
$ws.Range("D48").Value = "'2.43"
$ws.Range("E48").Value = "  +13.18%  "

# Row 49
$ws.Range("D49").Value = "'1.10"
$ws.Range("E49").Value = "  +3.43%  "

# Row 50
$ws.Range("D50").Value = "'1.13"
$ws.Range("E50").Value = "  +2.43%  "

# Row 51
$ws.Range("E51").Value = "  +1.65%  "
